# Update "last_edited_time" (column D) for the rows that were touched by
# this edit, and update the corresponding numeric metrics in row 13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LUY_KE_THANG_SOC_TRANG")

$newTimestamp = "2024-07-18T15:58:00.000Z"

# Rows whose last_edited_time (column D) changed.
$rows = @(2, 3, 6, 8, 11, 13)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 4).Value = $newTimestamp
}

# Row 13 numeric property updates ("cac bang" luy ke / thang 7 figures).
$ws.Range("S13").Value = 141172000
$ws.Range("W13").Value = 21416000
$ws.Range("AA13").Value = 8000000
$ws.Range("AE13").Value = 162588000
$ws.Range("AH13").Value = 137588000
$ws.Range("AK13").Value = 17
$ws.Range("AQ13").Value = 145588000
